# Re-simulated Week 17, factoring in more player injuries
# A.Thielen is removed from the Receiving sheet (e.g. ruled out with injury);
# the remaining receivers shift up to fill the gap.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Receiving")

# Remove A.Thielen's entire row (row 7), shifting everything below it up.
$ws.Rows.Item(7).Delete()

# Leave the workbook with the Receiving sheet active and I9 selected.
$ws.Activate()
$ws.Range("I9").Select()
